$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the locked cells below
# can be edited, then re-apply protection at the end so the sheet is left
# in the same (protected) state it started in.
$ws.Unprotect()

# Update the confidential/model-holdings date string in A11
# ("2021-05-25" -> "2021-05-26")
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# Update Weight (column D) / Percent Change (column E) values
$ws.Range("D2").Value = 0.4996237749781156
$ws.Range("E2").Value = 0.002084454007530878

$ws.Range("D3").Value = 0.2471074957538683
$ws.Range("E3").Value = 0.001742666279407556

$ws.Range("D4").Value = 0.09515744377716941
$ws.Range("E4").Value = 0.0068242133198535

$ws.Range("D5").Value = 0.1017869066057173
$ws.Range("E5").Value = 0.0111579934364745

$ws.Range("D6").Value = 0.02955971937545501
$ws.Range("E6").Value = 0.02111584657109655

$ws.Range("D7").Value = 0.02676465950967439
$ws.Range("E7").Value = 0.01945555817054978

$ws.Range("E8").Value = 0.004402080900405281

# Restore sheet protection to match the original workbook state.
$ws.Protect()
